$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet "1" -> "ქედა"
$ws.Name = "ქედა"

# Clear the subtitle row (row 2: "(მოსახლეობის აღწერის შედეგებით)") entirely, deleting it
$ws.Rows("2").Delete()

# Delete columns B and C (which held the 1989 / 2002 data), keeping only the 2014 column (shifts to column B)
$ws.Range("B1:C6").Delete()

$ws.Range("A2").Select() | Out-Null

$wb.Save()
